$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-18 18:29:30"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
